$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix data-entry mistakes where the Fotobond# value had been typed into
# column E ("Toevoeging Peter") instead of column F ("Fotobond#").
# Move the value over to column F and blank out column E.
$rowsToFix = 3, 49, 50, 53, 63
foreach ($r in $rowsToFix) {
    $val = $ws.Range("E$r").Value2
    $ws.Range("F$r").Value = $val
    $ws.Range("E$r").ClearContents()
}

# Row 27 had a stray "Nee" in the "Toevoeging Peter" column; clear it.
$ws.Range("E27").ClearContents()

# Filter the table on the "Fotobond#" column (column F) to show only the
# rows that are still marked "-" (i.e. not yet linked to a Fotobond number).
$lo = $ws.ListObjects.Item(1)
$lo.Range.AutoFilter(6, @("-"), 7)

# Update the selected cell shown in the sheet view.
$ws.Range("G27").Select()
